$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking price cells to retain their original text formatting
# (these columns store values like "1.006" / "41.48" as text, not numbers)
$textCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D21", "D22", "D24", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = '27.341.65'
$ws.Range("E2").Value = '  -1.74%  '
$ws.Range("D3").Value = '1.727.97'
$ws.Range("E3").Value = '  -2.05%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").Value = '321.17'
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("D7").Value = '0.4580'
$ws.Range("E7").Value = '  +7.54%  '
$ws.Range("D8").Value = '0.3517'
$ws.Range("E8").Value = '  -3.25%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '41.48'
$ws.Range("E9").Value = '  -2.47%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.07300'
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("D11").Value = '1.070'
$ws.Range("E11").Value = '  -2.19%  '
$ws.Range("D12").Value = '1.006'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").Value = '20.24'
$ws.Range("E13").Value = '  -2.48%  '
$ws.Range("D14").Value = '5.875'
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("D15").Value = '7.031'
$ws.Range("E15").Value = '  -3.50%  '
$ws.Range("D16").Value = '1.736.24'
$ws.Range("E16").Value = '  -1.79%  '
$ws.Range("D17").Value = '90.85'
$ws.Range("E17").Value = '  -0.72%  '
$ws.Range("D18").Value = '0.00001044'
$ws.Range("E18").Value = '  -1.83%  '
$ws.Range("D19").Value = '0.06329'
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").Value = '16.52'
$ws.Range("E21").Value = '  -3.14%  '
$ws.Range("D22").Value = '5.715'
$ws.Range("E22").Value = '  -3.35%  '
$ws.Range("D23").Value = '27.418.80'
$ws.Range("E23").Value = '  -1.60%  '
$ws.Range("D24").Value = '10.99'
$ws.Range("E24").Value = '  -2.21%  '
$ws.Range("E25").Value = '  -1.38%  '
$ws.Range("D26").Value = '162.00'
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("D27").Value = '19.74'
$ws.Range("E27").Value = '  -2.67%  '
$ws.Range("D28").Value = '1.931.01'
$ws.Range("E28").Value = '  -2.72%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '124.41'
$ws.Range("E29").Value = '  -0.85%  '
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = '2.032'
$ws.Range("E30").Value = '  -5.62%  '
$ws.Range("D31").Value = '1.033'
$ws.Range("E31").Value = '  -7.72%  '
$ws.Range("D32").Value = '0.09122'
$ws.Range("E32").Value = '  +2.63%  '
$ws.Range("D33").Value = '3.650'
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("D34").Value = '5.329'
$ws.Range("E34").Value = '  -4.28%  '
$ws.Range("D35").Value = '0.02253'
$ws.Range("E35").Value = '  -1.77%  '
$ws.Range("D36").Value = '11.55'
$ws.Range("E36").Value = '  -5.95%  '
$ws.Range("D37").Value = '0.05934'
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("D38").Value = '0.2037'
$ws.Range("E38").Value = '  -3.33%  '
$ws.Range("D39").Value = '0.6177'
$ws.Range("E39").Value = '  -2.67%  '
$ws.Range("D40").Value = '4.839'
$ws.Range("E40").Value = '  -2.76%  '
$ws.Range("D41").Value = '1.178'
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").Value = '1.370'
$ws.Range("E42").Value = '  -1.78%  '
$ws.Range("D43").Value = '7.652'
$ws.Range("E43").Value = '  -3.18%  '
$ws.Range("D44").Value = '13.04'
$ws.Range("E44").Value = '  -2.06%  '
$ws.Range("D45").Value = '3.680'
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("D46").Value = '0.5754'
$ws.Range("E46").Value = '  -2.20%  '
$ws.Range("D47").Value = '121.39'
$ws.Range("E47").Value = '  -1.31%  '
$ws.Range("D48").Value = '1.904'
$ws.Range("E48").Value = '  -4.22%  '
$ws.Range("D49").Value = '0.06820'
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("D50").Value = '1.102'
$ws.Range("E50").Value = '  -7.19%  '
$ws.Range("D51").Value = '70.70'
$ws.Range("E51").Value = '  -4.42%  '
